$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3798.3333
$ws.Cells.Item(64, 9).Value = 3798.5
$ws.Cells.Item(64, 10).Value = 3798
$ws.Cells.Item(64, 11).Value = 3798.5
$ws.Cells.Item(64, 12).Value = 3798
$ws.Cells.Item(64, 13).Value = -3550.5
$ws.Cells.Item(64, 14).Value = -4294
$ws.Cells.Item(67, 8).Value = 3798.3333
$ws.Cells.Item(67, 9).Value = 3798.5
$ws.Cells.Item(67, 10).Value = 3798
$ws.Cells.Item(67, 11).Value = 3798.5
$ws.Cells.Item(67, 12).Value = 3798
$ws.Cells.Item(67, 13).Value = -2940.5
$ws.Cells.Item(67, 14).Value = -5514
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(88, 8).Value = 1919.875
$ws.Cells.Item(88, 9).Value = 1935.5
$ws.Cells.Item(88, 11).Value = 1935.5
$ws.Cells.Item(88, 13).Value = -1529.5
$ws.Cells.Item(91, 8).Value = 1919.875
$ws.Cells.Item(91, 9).Value = 1935.5
$ws.Cells.Item(91, 11).Value = 1935.5
$ws.Cells.Item(91, 13).Value = -531.5
$ws.Cells.Item(98, 8).Value = 1295.6428
$ws.Cells.Item(98, 9).Value = 1348.9
$ws.Cells.Item(98, 10).Value = 1162.5
$ws.Cells.Item(98, 11).Value = 1348.9
$ws.Cells.Item(98, 12).Value = 1162.5
$ws.Cells.Item(98, 13).Value = 149.0999999999999
$ws.Cells.Item(98, 14).Value = -4158.5
$ws.Cells.Item(112, 8).Value = 3247
$ws.Cells.Item(112, 10).Value = 3321.8845
$ws.Cells.Item(112, 12).Value = 9965.6535
$ws.Cells.Item(112, 14).Value = -12181.6535
$ws.Cells.Item(122, 8).Value = 1295.6428
$ws.Cells.Item(122, 9).Value = 1348.9
$ws.Cells.Item(122, 10).Value = 1162.5
$ws.Cells.Item(122, 11).Value = 4046.7
$ws.Cells.Item(122, 12).Value = 3487.5
$ws.Cells.Item(122, 13).Value = -1596.7
$ws.Cells.Item(122, 14).Value = -8387.5
$ws.Cells.Item(125, 8).Value = 705.1667
$ws.Cells.Item(125, 10).Value = 779.6667
$ws.Cells.Item(125, 12).Value = 7017.0003
$ws.Cells.Item(125, 14).Value = -11937.0003
$ws.Cells.Item(127, 8).Value = 548.6667
$ws.Cells.Item(127, 9).Value = 458.4
$ws.Cells.Item(127, 11).Value = 1375.2
$ws.Cells.Item(127, 13).Value = 3584.8
$ws.Cells.Item(132, 8).Value = 911397
$ws.Cells.Item(132, 9).Value = 2536.7
$ws.Cells.Item(132, 10).Value = 10000000
$ws.Cells.Item(132, 11).Value = 7610.099999999999
$ws.Cells.Item(132, 12).Value = 30000000
$ws.Cells.Item(132, 13).Value = -5080.099999999999
$ws.Cells.Item(132, 14).Value = -30005060
$ws.Cells.Item(135, 8).Value = 2548.5
$ws.Cells.Item(135, 9).Value = 2548.5
$ws.Cells.Item(135, 11).Value = 22936.5
$ws.Cells.Item(135, 13).Value = -20401.5
$ws.Cells.Item(137, 8).Value = 3483.8965
$ws.Cells.Item(137, 9).Value = 3028.5
$ws.Cells.Item(137, 10).Value = 3723.5789
$ws.Cells.Item(137, 11).Value = 9085.5
$ws.Cells.Item(137, 12).Value = 11170.7367
$ws.Cells.Item(137, 13).Value = -6535.5
$ws.Cells.Item(137, 14).Value = -16270.7367
$ws.Cells.Item(138, 8).Value = 2889.6985
$ws.Cells.Item(138, 9).Value = 2204.75
$ws.Cells.Item(138, 10).Value = 3050.8628
$ws.Cells.Item(138, 11).Value = 6614.25
$ws.Cells.Item(138, 12).Value = 9152.588400000001
$ws.Cells.Item(138, 13).Value = -1474.25
$ws.Cells.Item(138, 14).Value = -19432.5884
$ws.Cells.Item(141, 8).Value = 4947.25
$ws.Cells.Item(141, 9).Value = 4485.222
$ws.Cells.Item(141, 10).Value = 6333.3335
$ws.Cells.Item(141, 11).Value = 13455.666
$ws.Cells.Item(141, 12).Value = 19000.0005
$ws.Cells.Item(141, 13).Value = -8275.665999999999
$ws.Cells.Item(141, 14).Value = -29360.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1619.826
$ws.Cells.Item(2, 9).Value = 1408.2632
$ws.Cells.Item(2, 11).Value = 1408.2632
$ws.Cells.Item(2, 13).Value = -1295.2632
$ws.Cells.Item(32, 8).Value = 12778.158
$ws.Cells.Item(32, 9).Value = 5687
$ws.Cells.Item(32, 11).Value = 5687
$ws.Cells.Item(32, 13).Value = -5400
$ws.Cells.Item(34, 8).Value = 11142212
$ws.Cells.Item(34, 10).Value = 39979
$ws.Cells.Item(34, 12).Value = 39979
$ws.Cells.Item(34, 14).Value = -40521
$ws.Cells.Item(97, 8).Value = 2682.9443
$ws.Cells.Item(97, 9).Value = 2689.2856
$ws.Cells.Item(97, 10).Value = 2660.75
$ws.Cells.Item(97, 11).Value = 2689.2856
$ws.Cells.Item(97, 12).Value = 2660.75
$ws.Cells.Item(97, 13).Value = -2193.2856
$ws.Cells.Item(97, 14).Value = -3652.75
$ws.Cells.Item(102, 8).Value = 2165.2
$ws.Cells.Item(102, 9).Value = 2165.2
$ws.Cells.Item(102, 11).Value = 2165.2
$ws.Cells.Item(102, 13).Value = -543.1999999999998
$ws.Cells.Item(116, 8).Value = 1619.826
$ws.Cells.Item(116, 9).Value = 1408.2632
$ws.Cells.Item(116, 11).Value = 1408.2632
$ws.Cells.Item(116, 13).Value = 885.7367999999999
$ws.Cells.Item(132, 8).Value = 1925.4
$ws.Cells.Item(132, 9).Value = 1657.125
$ws.Cells.Item(132, 11).Value = 4971.375
$ws.Cells.Item(132, 13).Value = -2441.375
$ws.Cells.Item(140, 8).Value = 76393
$ws.Cells.Item(140, 9).Value = 76393
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 76393
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -71213
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1619.826
$ws.Cells.Item(3, 9).Value = 1408.2632
$ws.Cells.Item(3, 11).Value = 1408.2632
$ws.Cells.Item(3, 13).Value = -1294.2632
$ws.Cells.Item(64, 8).Value = 1083.1818
$ws.Cells.Item(64, 10).Value = 1086.1666
$ws.Cells.Item(64, 12).Value = 1086.1666
$ws.Cells.Item(64, 14).Value = -1536.1666
$ws.Cells.Item(67, 8).Value = 1083.1818
$ws.Cells.Item(67, 10).Value = 1086.1666
$ws.Cells.Item(67, 12).Value = 1086.1666
$ws.Cells.Item(67, 14).Value = -2646.1666
$ws.Cells.Item(86, 8).Value = 2961.625
$ws.Cells.Item(86, 9).Value = 2961.625
$ws.Cells.Item(86, 11).Value = 2961.625
$ws.Cells.Item(86, 13).Value = -1838.625
$ws.Cells.Item(89, 8).Value = 2961.625
$ws.Cells.Item(89, 9).Value = 2961.625
$ws.Cells.Item(89, 11).Value = 14808.125
$ws.Cells.Item(89, 13).Value = -9192.125
$ws.Cells.Item(107, 8).Value = 1714.9231
$ws.Cells.Item(107, 9).Value = 1714.9231
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 1714.9231
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 205.0769
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 3861.8333
$ws.Cells.Item(134, 9).Value = 3834.3
$ws.Cells.Item(134, 11).Value = 11502.9
$ws.Cells.Item(134, 13).Value = -8967.900000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2663.3333
$ws.Cells.Item(58, 9).Value = 2764.3125
$ws.Cells.Item(58, 11).Value = 2764.3125
$ws.Cells.Item(58, 13).Value = -2561.3125
$ws.Cells.Item(96, 8).Value = 13140.667
$ws.Cells.Item(96, 10).Value = 13140.667
$ws.Cells.Item(96, 12).Value = 13140.667
$ws.Cells.Item(96, 14).Value = -18632.667
$ws.Cells.Item(99, 8).Value = 3362.0908
$ws.Cells.Item(99, 9).Value = 3770.5557
$ws.Cells.Item(99, 10).Value = 1524
$ws.Cells.Item(99, 11).Value = 3770.5557
$ws.Cells.Item(99, 12).Value = 1524
$ws.Cells.Item(99, 13).Value = -2272.5557
$ws.Cells.Item(99, 14).Value = -4520
$ws.Cells.Item(122, 8).Value = 2305.3635
$ws.Cells.Item(122, 9).Value = 2611.4707
$ws.Cells.Item(122, 10).Value = 1264.6
$ws.Cells.Item(122, 11).Value = 7834.4121
$ws.Cells.Item(122, 12).Value = 3793.8
$ws.Cells.Item(122, 13).Value = -5384.4121
$ws.Cells.Item(122, 14).Value = -8693.799999999999
$ws.Cells.Item(126, 8).Value = 3362.0908
$ws.Cells.Item(126, 9).Value = 3770.5557
$ws.Cells.Item(126, 10).Value = 1524
$ws.Cells.Item(126, 11).Value = 11311.6671
$ws.Cells.Item(126, 12).Value = 4572
$ws.Cells.Item(126, 13).Value = -8841.667099999999
$ws.Cells.Item(126, 14).Value = -9512
$ws.Cells.Item(132, 8).Value = 4796.3335
$ws.Cells.Item(132, 9).Value = 4796.3335
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 14389.0005
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -11859.0005
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 1458.9375
$ws.Cells.Item(134, 9).Value = 1612.8334
$ws.Cells.Item(134, 10).Value = 997.25
$ws.Cells.Item(134, 11).Value = 4838.5002
$ws.Cells.Item(134, 12).Value = 2991.75
$ws.Cells.Item(134, 13).Value = -2303.5002
$ws.Cells.Item(134, 14).Value = -8061.75
$ws.Cells.Item(136, 8).Value = 2663.3333
$ws.Cells.Item(136, 9).Value = 2764.3125
$ws.Cells.Item(136, 11).Value = 8292.9375
$ws.Cells.Item(136, 13).Value = -5742.9375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 97.25
$ws.Cells.Item(2, 9).Value = 96.333336
$ws.Cells.Item(2, 11).Value = 578.000016
$ws.Cells.Item(2, 13).Value = -465.000016
$ws.Cells.Item(8, 8).Value = 466.1111
$ws.Cells.Item(8, 9).Value = 466.1111
$ws.Cells.Item(8, 11).Value = 1398.3333
$ws.Cells.Item(8, 13).Value = -1259.3333
$ws.Cells.Item(12, 8).Value = 35.46154
$ws.Cells.Item(12, 9).Value = 34.333332
$ws.Cells.Item(12, 11).Value = 102.999996
$ws.Cells.Item(12, 13).Value = 70.000004
$ws.Cells.Item(23, 8).Value = 167.57143
$ws.Cells.Item(23, 9).Value = 154.8
$ws.Cells.Item(23, 11).Value = 464.4
$ws.Cells.Item(23, 13).Value = -229.4
$ws.Cells.Item(34, 8).Value = 740
$ws.Cells.Item(34, 10).Value = 935.8
$ws.Cells.Item(34, 12).Value = 2807.4
$ws.Cells.Item(34, 14).Value = -2975.4
$ws.Cells.Item(39, 8).Value = 3784.4285
$ws.Cells.Item(39, 10).Value = 5333
$ws.Cells.Item(39, 12).Value = 15999
$ws.Cells.Item(39, 14).Value = -16587
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(44, 8).Value = 6746.6665
$ws.Cells.Item(44, 9).Value = 6492.5
$ws.Cells.Item(44, 11).Value = 19477.5
$ws.Cells.Item(44, 13).Value = -19079.5
$ws.Cells.Item(55, 8).Value = 224414.89
$ws.Cells.Item(55, 9).Value = 559
$ws.Cells.Item(55, 10).Value = 403499.6
$ws.Cells.Item(55, 11).Value = 1677
$ws.Cells.Item(55, 12).Value = 1210498.8
$ws.Cells.Item(55, 13).Value = -1500
$ws.Cells.Item(55, 14).Value = -1210852.8
$ws.Cells.Item(62, 8).Value = 6528.7144
$ws.Cells.Item(62, 9).Value = 1856
$ws.Cells.Item(62, 10).Value = 8397.799999999999
$ws.Cells.Item(62, 11).Value = 5568
$ws.Cells.Item(62, 12).Value = 25193.4
$ws.Cells.Item(62, 13).Value = -4882
$ws.Cells.Item(62, 14).Value = -26565.4
$ws.Cells.Item(63, 8).Value = 7999.75
$ws.Cells.Item(63, 10).Value = 10000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 14).Value = -31498
$ws.Cells.Item(65, 8).Value = 6528.7144
$ws.Cells.Item(65, 9).Value = 1856
$ws.Cells.Item(65, 10).Value = 8397.799999999999
$ws.Cells.Item(65, 11).Value = 16704
$ws.Cells.Item(65, 12).Value = 75580.2
$ws.Cells.Item(65, 13).Value = -13272
$ws.Cells.Item(65, 14).Value = -82444.2
$ws.Cells.Item(66, 8).Value = 7999.75
$ws.Cells.Item(66, 10).Value = 10000
$ws.Cells.Item(66, 12).Value = 90000
$ws.Cells.Item(66, 14).Value = -97488
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 407.91666
$ws.Cells.Item(107, 9).Value = 264.75
$ws.Cells.Item(107, 10).Value = 479.5
$ws.Cells.Item(107, 11).Value = 794.25
$ws.Cells.Item(107, 12).Value = 1438.5
$ws.Cells.Item(107, 13).Value = 1125.75
$ws.Cells.Item(107, 14).Value = -5278.5
$ws.Cells.Item(122, 8).Value = 2550.7778
$ws.Cells.Item(122, 9).Value = 2349
$ws.Cells.Item(122, 10).Value = 2712.2
$ws.Cells.Item(122, 11).Value = 21141
$ws.Cells.Item(122, 12).Value = 24409.8
$ws.Cells.Item(122, 13).Value = -18691
$ws.Cells.Item(122, 14).Value = -29309.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).ClearContents()
$ws.Cells.Item(97, 8).Value = 1828.5834
$ws.Cells.Item(97, 9).Value = 1828.5834
$ws.Cells.Item(97, 11).Value = 1828.5834
$ws.Cells.Item(97, 13).Value = -1332.5834
$ws.Cells.Item(102, 8).Value = 3340.48
$ws.Cells.Item(102, 9).Value = 2784.1765
$ws.Cells.Item(102, 10).Value = 4522.625
$ws.Cells.Item(102, 11).Value = 2784.1765
$ws.Cells.Item(102, 12).Value = 4522.625
$ws.Cells.Item(102, 13).Value = -1162.1765
$ws.Cells.Item(102, 14).Value = -7766.625
$ws.Cells.Item(113, 8).Value = 3818.4285
$ws.Cells.Item(113, 9).Value = 2822
$ws.Cells.Item(113, 11).Value = 2822
$ws.Cells.Item(113, 13).Value = -652
$ws.Cells.Item(122, 8).Value = 41734.6
$ws.Cells.Item(122, 9).Value = 126176.5
$ws.Cells.Item(122, 11).Value = 378529.5
$ws.Cells.Item(122, 13).Value = -376079.5
$ws.Cells.Item(132, 8).Value = 3395.4644
$ws.Cells.Item(132, 9).Value = 4122.6665
$ws.Cells.Item(132, 10).Value = 2086.5
$ws.Cells.Item(132, 11).Value = 12367.9995
$ws.Cells.Item(132, 12).Value = 6259.5
$ws.Cells.Item(132, 13).Value = -9837.999500000002
$ws.Cells.Item(132, 14).Value = -11319.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7268.385
$ws.Cells.Item(7, 9).Value = 6070.1
$ws.Cells.Item(7, 11).Value = 6070.1
$ws.Cells.Item(7, 13).Value = -5958.1
$ws.Cells.Item(22, 8).Value = 811.2
$ws.Cells.Item(22, 9).Value = 764.25
$ws.Cells.Item(22, 10).Value = 999
$ws.Cells.Item(22, 11).Value = 764.25
$ws.Cells.Item(22, 12).Value = 999
$ws.Cells.Item(22, 13).Value = -469.25
$ws.Cells.Item(22, 14).Value = -1589
$ws.Cells.Item(26, 8).Value = 7272.727
$ws.Cells.Item(26, 9).Value = 6750
$ws.Cells.Item(26, 10).Value = 8666.666999999999
$ws.Cells.Item(26, 11).Value = 6750
$ws.Cells.Item(26, 12).Value = 8666.666999999999
$ws.Cells.Item(26, 13).Value = -6455
$ws.Cells.Item(26, 14).Value = -9256.666999999999
$ws.Cells.Item(27, 8).Value = 811.2
$ws.Cells.Item(27, 9).Value = 764.25
$ws.Cells.Item(27, 10).Value = 999
$ws.Cells.Item(27, 11).Value = 764.25
$ws.Cells.Item(27, 12).Value = 999
$ws.Cells.Item(27, 13).Value = -657.25
$ws.Cells.Item(27, 14).Value = -1213
$ws.Cells.Item(40, 8).Value = 4390.1816
$ws.Cells.Item(40, 9).Value = 3411.5
$ws.Cells.Item(40, 10).Value = 7000
$ws.Cells.Item(40, 11).Value = 3411.5
$ws.Cells.Item(40, 12).Value = 7000
$ws.Cells.Item(40, 13).Value = -3275.5
$ws.Cells.Item(40, 14).Value = -7272
$ws.Cells.Item(68, 8).Value = 2261.5
$ws.Cells.Item(68, 9).Value = 2156.2856
$ws.Cells.Item(68, 11).Value = 2156.2856
$ws.Cells.Item(68, 13).Value = -1407.2856
$ws.Cells.Item(71, 8).Value = 2261.5
$ws.Cells.Item(71, 9).Value = 2156.2856
$ws.Cells.Item(71, 11).Value = 10781.428
$ws.Cells.Item(71, 13).Value = -7037.428
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(93, 8).Value = 3009.1
$ws.Cells.Item(93, 9).Value = 3024.9167
$ws.Cells.Item(93, 10).Value = 2985.375
$ws.Cells.Item(93, 11).Value = 3024.9167
$ws.Cells.Item(93, 12).Value = 2985.375
$ws.Cells.Item(93, 13).Value = -1776.9167
$ws.Cells.Item(93, 14).Value = -5481.375
$ws.Cells.Item(100, 8).Value = 3363.6843
$ws.Cells.Item(100, 9).Value = 3346.8462
$ws.Cells.Item(100, 10).Value = 3400.1667
$ws.Cells.Item(100, 11).Value = 3346.8462
$ws.Cells.Item(100, 12).Value = 3400.1667
$ws.Cells.Item(100, 13).Value = -2805.8462
$ws.Cells.Item(100, 14).Value = -4482.1667
$ws.Cells.Item(126, 8).Value = 7268.385
$ws.Cells.Item(126, 9).Value = 6070.1
$ws.Cells.Item(126, 11).Value = 18210.3
$ws.Cells.Item(126, 13).Value = -15740.3
$ws.Cells.Item(136, 8).Value = 4269.522
$ws.Cells.Item(136, 9).Value = 4177.35
$ws.Cells.Item(136, 11).Value = 12532.05
$ws.Cells.Item(136, 13).Value = -9982.050000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).ClearContents()
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).ClearContents()
$ws.Cells.Item(84, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 18698.436
$ws.Cells.Item(122, 9).Value = 17669.809
$ws.Cells.Item(122, 11).Value = 53009.427
$ws.Cells.Item(122, 13).Value = -50559.427
$ws.Cells.Item(126, 8).Value = 2784.5
$ws.Cells.Item(126, 9).Value = 2359.0588
$ws.Cells.Item(126, 10).Value = 5195.3335
$ws.Cells.Item(126, 11).Value = 7077.176399999999
$ws.Cells.Item(126, 12).Value = 15586.0005
$ws.Cells.Item(126, 13).Value = -4607.176399999999
$ws.Cells.Item(126, 14).Value = -20526.0005
$ws.Cells.Item(132, 8).Value = 2679.15
$ws.Cells.Item(132, 9).Value = 2557.0527
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 7671.158100000001
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -5141.158100000001
$ws.Cells.Item(132, 14).Value = -20057
$ws.Cells.Item(136, 8).Value = 3153.5715
$ws.Cells.Item(136, 9).Value = 3083.8333
$ws.Cells.Item(136, 11).Value = 9251.499899999999
$ws.Cells.Item(136, 13).Value = -6701.499899999999
